$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.7594023393015787, 0.4843090051523787)
    3  = @(0.9191225853813544, 1.36548478079573)
    4  = @(0.7020093069832635, 0.2683431762503036)
    5  = @(0.8321536731438657, 0.7941101279597217)
    6  = @(0.6528774572507436, 0.393751305815818)
    7  = @(0.7224392769862428, 0.3260101172571158)
    8  = @(0.6689448310006813, 0.2490139508291321)
    9  = @(0.6495829054855937, 0.29563841016148)
    10 = @(0.8572950184359456, 1.08648801352677)
    11 = @(1.166168156856927, 1.193618170623214)
    12 = @(1.749765316945811, 1.778185541566537)
    13 = @(1.481843511032632, 1.545785333250495)
    14 = @(0.633761718521524, 0.7354779328243995)
    15 = @(0.6813228819063484, 0.3941211788185811)
    16 = @(0.8671805920593697, 0.5340462395551343)
    17 = @(0.7478218636558492, 0.2504249277958301)
    18 = @(0.7037996867421685, 0.3114347773329301)
    19 = @(0.6588107456362138, 0.5303400442387239)
    20 = @(0.7430843995336209, 0.2518973143680697)
    21 = @(0.769638303072412, 0.2629558359801659)
    22 = @(0.8159760480750113, 0.3591514826768604)
    23 = @(0.8037184033040029, 0.3251386563633973)
    24 = @(0.6400911217032833, 0.7399460337843894)
    25 = @(0.7676949679739734, 0.260684823520987)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("K$row").Value = $pair[0]
    $ws.Range("M$row").Value = $pair[1]
}
